# Scheduled market-price refresh: update currentAveragePrice / Leve profit
# columns (H:N) for the affected leve rows across the ALC, ARM, CRP, GSM and
# LTW sheets. One row (ALC!17) also drops its now-unused LeveProfitNQ (M)
# value entirely, matching the upstream price-feed output for that row.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 53659.316
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 53659.316
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 160977.948
$ws.Range("N17").Value = -161313.948
$ws.Range("H28").Value = 265.38095
$ws.Range("I28").Value = 258.72223
$ws.Range("J28").Value = 305.33334
$ws.Range("K28").Value = 258.72223
$ws.Range("L28").Value = 305.33334
$ws.Range("M28").Value = 226.27777
$ws.Range("N28").Value = -1275.33334
$ws.Range("H86").Value = 4372718.5
$ws.Range("I86").Value = 51500
$ws.Range("J86").Value = 10134343
$ws.Range("K86").Value = 51500
$ws.Range("L86").Value = 10134343
$ws.Range("M86").Value = -50377
$ws.Range("N86").Value = -10136589
$ws.Range("H89").Value = 4372718.5
$ws.Range("I89").Value = 51500
$ws.Range("J89").Value = 10134343
$ws.Range("K89").Value = 257500
$ws.Range("L89").Value = 50671715
$ws.Range("M89").Value = -251884
$ws.Range("N89").Value = -50682947
$ws.Range("H98").Value = 125000700
$ws.Range("I98").Value = 156250510
$ws.Range("J98").Value = 1500
$ws.Range("K98").Value = 156250510
$ws.Range("L98").Value = 1500
$ws.Range("M98").Value = -156249012
$ws.Range("N98").Value = -4496
$ws.Range("H106").Value = 200002690
$ws.Range("I106").Value = 1000000000
$ws.Range("J106").Value = 3350
$ws.Range("K106").Value = 1000000000
$ws.Range("L106").Value = 3350
$ws.Range("M106").Value = -999999369
$ws.Range("N106").Value = -4612
$ws.Range("H121").Value = 799.1111
$ws.Range("I121").Value = 309.25
$ws.Range("J121").Value = 1191
$ws.Range("K121").Value = 927.75
$ws.Range("L121").Value = 3573
$ws.Range("M121").Value = 819.25
$ws.Range("N121").Value = -7067
$ws.Range("H122").Value = 125000700
$ws.Range("I122").Value = 156250510
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 468751530
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -468749080
$ws.Range("N122").Value = -9400
$ws.Range("H135").Value = 717.5472
$ws.Range("I135").Value = 732.7083
$ws.Range("J135").Value = 572
$ws.Range("K135").Value = 6594.3747
$ws.Range("L135").Value = 5148
$ws.Range("M135").Value = -4059.3747
$ws.Range("N135").Value = -10218
$ws.Range("H137").Value = 1228.4681
$ws.Range("I137").Value = 1040.4651
$ws.Range("K137").Value = 3121.3953
$ws.Range("M137").Value = -571.3952999999997
$ws.Range("H138").Value = 2058.4138
$ws.Range("I138").Value = 1203.2703
$ws.Range("J138").Value = 3565.0952
$ws.Range("K138").Value = 3609.810899999999
$ws.Range("L138").Value = 10695.2856
$ws.Range("M138").Value = 1530.189100000001
$ws.Range("N138").Value = -20975.2856
$ws.Range("M17").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12630814
$ws.Range("I32").Value = 1505677.9
$ws.Range("J32").Value = 71435100
$ws.Range("K32").Value = 1505677.9
$ws.Range("L32").Value = 71435100
$ws.Range("M32").Value = -1505390.9
$ws.Range("N32").Value = -71435674
$ws.Range("H61").Value = 1366.7106
$ws.Range("I61").Value = 1040.931
$ws.Range("J61").Value = 2416.4443
$ws.Range("K61").Value = 1040.931
$ws.Range("L61").Value = 2416.4443
$ws.Range("M61").Value = -828.931
$ws.Range("N61").Value = -2840.4443
$ws.Range("H132").Value = 22529688
$ws.Range("I132").Value = 24391166
$ws.Range("K132").Value = 73173498
$ws.Range("M132").Value = -73170968
$ws.Range("H136").Value = 1366.7106
$ws.Range("I136").Value = 1040.931
$ws.Range("J136").Value = 2416.4443
$ws.Range("K136").Value = 3122.793
$ws.Range("L136").Value = 7249.3329
$ws.Range("M136").Value = -572.7930000000001
$ws.Range("N136").Value = -12349.3329

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 14354.4
$ws.Range("I105").Value = 21713.889
$ws.Range("K105").Value = 21713.889
$ws.Range("M105").Value = -19966.889
$ws.Range("H132").Value = 5209240
$ws.Range("I132").Value = 744.4
$ws.Range("J132").Value = 23811010
$ws.Range("K132").Value = 2233.2
$ws.Range("L132").Value = 71433030
$ws.Range("M132").Value = 296.8000000000002
$ws.Range("N132").Value = -71438090

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4910.3477
$ws.Range("I70").Value = 4075.5625
$ws.Range("K70").Value = 4075.5625
$ws.Range("M70").Value = -3805.5625
$ws.Range("H73").Value = 4910.3477
$ws.Range("I73").Value = 4075.5625
$ws.Range("K73").Value = 4075.5625
$ws.Range("M73").Value = -3139.5625
$ws.Range("H93").Value = 20251
$ws.Range("J93").Value = 20251
$ws.Range("L93").Value = 20251
$ws.Range("N93").Value = -23995
$ws.Range("H113").Value = 1216.6666
$ws.Range("I113").Value = 1100
$ws.Range("J113").Value = 1275
$ws.Range("K113").Value = 1100
$ws.Range("L113").Value = 1275
$ws.Range("M113").Value = 1070
$ws.Range("N113").Value = -5615
$ws.Range("H122").Value = 25005934
$ws.Range("I122").Value = 29418454
$ws.Range("J122").Value = 1650
$ws.Range("K122").Value = 88255362
$ws.Range("L122").Value = 4950
$ws.Range("M122").Value = -88252912
$ws.Range("N122").Value = -9850
$ws.Range("H132").Value = 7502
$ws.Range("I132").Value = 5328.273
$ws.Range("J132").Value = 12284.2
$ws.Range("K132").Value = 15984.819
$ws.Range("L132").Value = 36852.60000000001
$ws.Range("M132").Value = -13454.819
$ws.Range("N132").Value = -41912.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 9617038
$ws.Range("I40").Value = 1730.7
$ws.Range("J40").Value = 41668064
$ws.Range("K40").Value = 1730.7
$ws.Range("L40").Value = 41668064
$ws.Range("M40").Value = -1594.7
$ws.Range("N40").Value = -41668336
$ws.Range("H69").Value = 33592
$ws.Range("J69").Value = 33592
$ws.Range("L69").Value = 33592
$ws.Range("N69").Value = -35214
$ws.Range("H72").Value = 33592
$ws.Range("J72").Value = 33592
$ws.Range("L72").Value = 100776
$ws.Range("N72").Value = -108888
$ws.Range("H132").Value = 18873760
$ws.Range("I132").Value = 27779322
$ws.Range("J132").Value = 14923.706
$ws.Range("K132").Value = 83337966
$ws.Range("L132").Value = 44771.118
$ws.Range("M132").Value = -83335436
$ws.Range("N132").Value = -49831.118
